$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.114.27'
$ws.Range("E2").Value = '  -1.65%  '
$ws.Range("D3").Value = '2.254.44'
$ws.Range("E3").Value = '  -3.46%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '298.09'
$ws.Range("E5").Value = '  -2.74%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '94.02'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.496'
$ws.Range("E7").Value = '  -2.68%  '
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.492'
$ws.Range("E9").Value = '  -3.51%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.75'
$ws.Range("E10").Value = '  -6.40%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0783'
$ws.Range("E11").Value = '  -2.06%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '48.25'
$ws.Range("E12").Value = '  -7.34%  '
$ws.Range("E13").Value = '  +0.54%  '
$ws.Range("E14").Value = '  -2.89%  '
$ws.Range("D15").Value = '2.602.93'
$ws.Range("E15").Value = '  -3.52%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.30'
$ws.Range("E16").Value = '  -3.18%  '
$ws.Range("D17").Value = '2.252.64'
$ws.Range("E17").Value = '  -1.04%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.771'
$ws.Range("E18").Value = '  -3.09%  '
$ws.Range("D19").Value = '42.041.27'
$ws.Range("E19").Value = '  -1.58%  '
$ws.Range("D20").Value = '0.0₃0887'
$ws.Range("E20").Value = '  -2.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.98'
$ws.Range("E22").Value = '  -4.59%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.25'
$ws.Range("E23").Value = '  -2.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '232.77'
$ws.Range("E24").Value = '  -1.63%  '
$ws.Range("E25").Value = '  -4.32%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.46'
$ws.Range("E27").Value = '  -4.25%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.75'
$ws.Range("E28").Value = '  -4.90%  '
$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '165.80'
$ws.Range("E29").Value = '  +3.97%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.04'
$ws.Range("E30").Value = '  -6.80%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '33.58'
$ws.Range("E31").Value = '  -3.89%  '
$ws.Range("E32").Value = '  -3.63%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  +0.04%  '
$ws.Range("E34").Value = '  -4.01%  '
$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.32'
$ws.Range("E35").Value = '  -5.20%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0690'
$ws.Range("E36").Value = '  -5.20%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.32'
$ws.Range("E37").Value = '  -6.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.79'
$ws.Range("E38").Value = '  -5.78%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.93'
$ws.Range("E39").Value = '  -8.54%  '
$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.109'
$ws.Range("E40").Value = '  -2.89%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0978'
$ws.Range("E41").Value = '  -5.17%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.71'
$ws.Range("E42").Value = '  -8.30%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.43'
$ws.Range("E43").Value = '  +3.19%  '
$ws.Range("D44").Value = '1.936.70'
$ws.Range("E44").Value = '  -3.90%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0277'
$ws.Range("E45").Value = '  -2.42%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.26'
$ws.Range("E46").Value = '  -8.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.49'
$ws.Range("E47").Value = '  -7.97%  '
$ws.Range("E49").Value = '  -3.29%  '
$ws.Range("D50").Value = '2.480.91'
$ws.Range("E50").Value = '  -3.04%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '51.91'
$ws.Range("E51").Value = '  -7.49%  '
